$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.417.86'
$ws.Cells.Item(2, 5).Value = '  -1.06%  '
$ws.Cells.Item(3, 4).Value = '2.052.07'
$ws.Cells.Item(3, 5).Value = '  -1.60%  '
$ws.Cells.Item(4, 4).Value = '''0.997'
$ws.Cells.Item(4, 5).Value = '  -0.54%  '
$ws.Cells.Item(5, 4).Value = '''230.28'
$ws.Cells.Item(5, 5).Value = '  -1.48%  '
$ws.Cells.Item(6, 4).Value = '''0.613'
$ws.Cells.Item(6, 5).Value = '  -1.86%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).Value = '''57.09'
$ws.Cells.Item(8, 5).Value = '  -2.58%  '
$ws.Cells.Item(9, 4).Value = '''0.386'
$ws.Cells.Item(9, 5).Value = '  -2.13%  '
$ws.Cells.Item(10, 4).Value = '''0.0809'
$ws.Cells.Item(10, 5).Value = '  +3.11%  '
$ws.Cells.Item(11, 5).Value = '  -1.84%  '
$ws.Cells.Item(12, 4).Value = '''14.69'
$ws.Cells.Item(12, 5).Value = '  -2.73%  '
$ws.Cells.Item(13, 4).Value = '2.356.32'
$ws.Cells.Item(13, 5).Value = '  -1.30%  '
$ws.Cells.Item(14, 4).Value = '''20.79'
$ws.Cells.Item(14, 5).Value = '  -2.44%  '
$ws.Cells.Item(15, 4).Value = '''0.758'
$ws.Cells.Item(15, 5).Value = '  -2.96%  '
$ws.Cells.Item(16, 4).Value = '''5.31'
$ws.Cells.Item(16, 5).Value = '  -1.36%  '
$ws.Cells.Item(17, 4).Value = '2.053.59'
$ws.Cells.Item(17, 5).Value = '  -1.33%  '
$ws.Cells.Item(18, 4).Value = '37.302.30'
$ws.Cells.Item(18, 5).Value = '  -1.18%  '
$ws.Cells.Item(19, 4).Value = '''6.04'
$ws.Cells.Item(19, 5).Value = '  -1.65%  '
$ws.Cells.Item(20, 4).Value = '''69.99'
$ws.Cells.Item(20, 5).Value = '  -1.79%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0836'
$ws.Cells.Item(21, 5).Value = '  -0.22%  '
$ws.Cells.Item(22, 4).Value = '''227.14'
$ws.Cells.Item(22, 5).Value = '  -1.47%  '
$ws.Cells.Item(23, 4).Value = '''0.999'
$ws.Cells.Item(23, 5).Value = '  +0.14%  '
$ws.Cells.Item(24, 4).Value = '''2.36'
$ws.Cells.Item(24, 5).Value = '  -1.32%  '
$ws.Cells.Item(25, 4).Value = '''2.29'
$ws.Cells.Item(25, 5).Value = '  -4.82%  '
$ws.Cells.Item(26, 4).Value = '''9.53'
$ws.Cells.Item(26, 5).Value = '  -3.01%  '
$ws.Cells.Item(27, 4).Value = '''167.84'
$ws.Cells.Item(27, 5).Value = '  -2.45%  '
$ws.Cells.Item(28, 4).Value = '''1.40'
$ws.Cells.Item(28, 5).Value = '  -0.36%  '
$ws.Cells.Item(29, 5).Value = '  -5.81%  '
$ws.Cells.Item(30, 4).Value = '''18.99'
$ws.Cells.Item(30, 5).Value = '  -2.82%  '
$ws.Cells.Item(31, 5).Value = '  -2.77%  '
$ws.Cells.Item(32, 4).Value = '''4.55'
$ws.Cells.Item(32, 5).Value = '  -3.78%  '
$ws.Cells.Item(33, 5).Value = '  -1.75%  '
$ws.Cells.Item(34, 4).Value = '''0.0615'
$ws.Cells.Item(34, 5).Value = '  -3.07%  '
$ws.Cells.Item(35, 4).Value = '''2.41'
$ws.Cells.Item(35, 5).Value = '  -2.83%  '
$ws.Cells.Item(36, 5).Value = '  -0.06%  '
$ws.Cells.Item(37, 4).Value = '''0.994'
$ws.Cells.Item(37, 5).Value = '  -0.60%  '
$ws.Cells.Item(38, 4).Value = '''3.24'
$ws.Cells.Item(38, 5).Value = '  -4.77%  '
$ws.Cells.Item(39, 5).Value = '  -0.46%  '
$ws.Cells.Item(40, 5).Value = '  -6.50%  '
$ws.Cells.Item(41, 4).Value = '''17.07'
$ws.Cells.Item(41, 5).Value = '  +1.25%  '
$ws.Cells.Item(42, 4).Value = '1.488.26'
$ws.Cells.Item(42, 5).Value = '  +1.91%  '
$ws.Cells.Item(43, 5).Value = '  -1.27%  '
$ws.Cells.Item(44, 5).Value = '  -3.26%  '
$ws.Cells.Item(45, 4).Value = '''96.68'
$ws.Cells.Item(45, 5).Value = '  -6.00%  '
$ws.Cells.Item(46, 5).Value = '  +0.87%  '
$ws.Cells.Item(47, 5).Value = '  -4.17%  '
$ws.Cells.Item(48, 5).Value = '  -2.38%  '
$ws.Cells.Item(49, 5).Value = '  -2.59%  '
$ws.Cells.Item(50, 4).Value = '''3.75'
$ws.Cells.Item(50, 5).Value = '  -8.82%  '
$ws.Cells.Item(51, 4).Value = '2.243.51'
$ws.Cells.Item(51, 5).Value = '  -1.36%  '
